$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update capacity column (D) values from 100 to 200 for rows 2-7
$ws.Range("D2:D7").Value = 200

# Update the active cell selection to D8
$ws.Range("D8").Select()
